$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format while we write values, so Excel does not
# auto-convert numeric-looking strings (e.g. "1.004", "314.30") into numbers
# and strip meaningful trailing zeros / dot-grouping. Reset the style back to
# "Normal" afterwards so the saved cells carry no explicit style (matching the
# original workbook formatting).
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "24.591.39"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "1.689.86"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "314.30"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").Value = "0.3887"
$ws.Range("E7").Value = "  -1.75%  "
$ws.Range("D8").Value = "0.4034"
$ws.Range("E8").Value = "  -0.86%  "
$ws.Range("D9").Value = "1.493"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("D11").Value = "53.05"
$ws.Range("E11").Value = "  +1.89%  "
$ws.Range("D12").Value = "0.08740"
$ws.Range("E12").Value = "  -1.77%  "
$ws.Range("E13").Value = "  +7.02%  "
$ws.Range("D14").Value = "7.520"
$ws.Range("E14").Value = "  +3.50%  "
$ws.Range("D15").Value = "0.00001357"
$ws.Range("D16").Value = "7.941"
$ws.Range("E16").Value = "  -1.87%  "
$ws.Range("D17").Value = "1.693.62"
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("D19").Value = "0.07101"
$ws.Range("E19").Value = "  +0.92%  "
$ws.Range("D20").Value = "19.94"
$ws.Range("E20").Value = "  +1.56%  "
$ws.Range("D21").Value = "7.257"
$ws.Range("E21").Value = "  +3.49%  "
$ws.Range("E22").Value = "  -0.40%  "
$ws.Range("D23").Value = "14.24"
$ws.Range("E23").Value = "  -0.86%  "
$ws.Range("D24").Value = "24.588.64"
$ws.Range("E24").Value = "  -0.41%  "
$ws.Range("D25").Value = "2.976"
$ws.Range("D26").Value = "2.355"
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("D27").Value = "22.71"
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").Value = "161.85"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("D29").Value = "8.760"
$ws.Range("E29").Value = "  +15.87%  "
$ws.Range("D30").Value = "136.69"
$ws.Range("E30").Value = "  +0.54%  "
$ws.Range("D31").Value = "5.222"
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("D32").Value = "1.878.09"
$ws.Range("E32").Value = "  -0.34%  "
$ws.Range("D33").Value = "0.08820"
$ws.Range("E33").Value = "  +2.19%  "
$ws.Range("D34").Value = "7.374"
$ws.Range("E34").Value = "  +4.36%  "
$ws.Range("D35").Value = "1.033"
$ws.Range("E35").Value = "  -2.20%  "
$ws.Range("D36").Value = "1.958"
$ws.Range("E36").Value = "  +3.87%  "
$ws.Range("D37").Value = "0.2747"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("D38").Value = "0.02910"
$ws.Range("E38").Value = "  +6.71%  "
$ws.Range("D39").Value = "10.76"
$ws.Range("E39").Value = "  -5.59%  "
$ws.Range("D40").Value = "14.22"
$ws.Range("E40").Value = "  -1.90%  "
$ws.Range("D41").Value = "0.09119"
$ws.Range("E41").Value = "  -1.21%  "
$ws.Range("D42").Value = "0.7895"
$ws.Range("E42").Value = "  +2.77%  "
$ws.Range("D43").Value = "1.457"
$ws.Range("E43").Value = "  -1.03%  "
$ws.Range("D44").Value = "16.75"
$ws.Range("E44").Value = "  +3.71%  "
$ws.Range("D45").Value = "0.7201"
$ws.Range("E45").Value = "  +0.46%  "
$ws.Range("D46").Value = "2.588"
$ws.Range("E46").Value = "  -0.64%  "
$ws.Range("D47").Value = "4.201"
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("D49").Value = "1.338"
$ws.Range("E49").Value = "  +1.27%  "
$ws.Range("D50").Value = "137.90"
$ws.Range("E50").Value = "  -1.85%  "
$ws.Range("D51").Value = "91.10"
$ws.Range("E51").Value = "  +0.14%  "

# Restore default ("Normal") style on column D so no stray cell-level style
# survives the text-format trick above.
$priceRange.Style = "Normal"

